$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 "I0" and J1 "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-15
$values = @{
    2  = @(1, 3)
    3  = @(2, 5)
    4  = @(2, 6)
    5  = @(3, 7)
    6  = @(1, 4)
    7  = @(2, 2)
    8  = @(1, 4)
    9  = @(1, 4)
    10 = @(3, 5)
    11 = @(6, 7)
    12 = @(8, 9)
    13 = @(8, 8)
    14 = @(6, 7)
    15 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
